$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both list the same rows for
# rows 5, 6, 10, 11, but the row offsets for the later block differ by 2
# between the two sheets (36/37/38 vs 38/39/40) because "全部类型" has a
# couple of extra rows earlier in the sheet.

$updates1 = @{
    5  = 1327
    6  = 18527
    10 = 6938
    11 = 357
    36 = 12186
    37 = 1297
    38 = 18
}

$updates4 = @{
    5  = 1327
    6  = 18527
    10 = 6938
    11 = 357
    38 = 12186
    39 = 1297
    40 = 18
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
